$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.499.74'
$ws.Range("E2").Value = '  -1.30%  '

$ws.Range("D3").Value = '3.082.16'
$ws.Range("E3").Value = '  -2.26%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '589.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.02'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.01%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.539'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.63%  '

$ws.Range("D9").Value = '3.072.50'
$ws.Range("E9").Value = '  -2.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.81'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.77%  '

$ws.Range("E12").Value = '  -1.38%  '

$ws.Range("B13").Value = 'ShibaInu'
$ws.Range("C13").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000240'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.23%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.47%  '

$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '3.591.11'
$ws.Range("E15").Value = '  -2.27%  '

$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.119'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.14'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.04%  '

$ws.Range("D18").Value = '63.496.32'
$ws.Range("E18").Value = '  -1.04%  '

$ws.Range("D19").Value = '3.075.47'
$ws.Range("E19").Value = '  -2.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '474.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.15%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.712'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.32%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.49%  '

$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.51%  '

$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '80.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.66%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.64%  '

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.62%  '

$ws.Range("E30").Value = '  -0.12%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.67'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.115'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.17'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.16%  '

$ws.Range("D35").Value = '0.0₃0850'
$ws.Range("E35").Value = '  -2.84%  '

$ws.Range("B36").Value = 'dogwifhat'
$ws.Range("C36").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.44'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.51%  '

$ws.Range("B37").Value = 'Mantle'
$ws.Range("C37").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.05'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.18%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.05'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.85%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.21'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.45%  '

$ws.Range("B40").Value = 'Cosmos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.27%  '

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.61'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '448.60'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.81%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.287'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.03%  '

$ws.Range("B44").Value = 'Arweave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.09%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0361'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.12%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.111'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.15%  '

$ws.Range("D47").Value = '2.805.33'
$ws.Range("E47").Value = '  -3.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.92'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.07%  '

$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.06%  '
